# This change reorders the data rows (rows 2-23, columns A-F) of the sheet.
# The underlying bug fix ("fixed a bug in mergePositionCollection") resulted
# in the position rows being merged/sorted in a different order; the data
# itself (each row keyed by its column-A "symbol" id) is unchanged, only the
# row order differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(1202, 2, 10, 10, 10, 10),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(801, 3, 67, 65, 52, 45),
    @(101, 9, 30, 15, 60, 15),
    @(701, 3, 90, 45, 97, 15),
    @(501, 9, 52, 30, 75, 45),
    @(401, 9, 48, 67, 75, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(301, 6, 45, 30, 60, 45),
    @(502, 0, 4, 0, 0, 0),
    @(2, 0, 2, 2, 2, 2),
    @(1, 0, 2, 2, 2, 2),
    @(1101, 0, 15, 30, 30, 0),
    @(802, 0, 4, 5, 4, 0),
    @(3, 0, 3, 3, 3, 3),
    @(402, 0, 0, 4, 0, 0),
    @(602, 0, 0, 4, 0, 9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Count; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $values[$col]
    }
}
